# watchlist new test cases implementations
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- Existing rows: Runmode column (C) changes from Y to N ---
$ws.Range("C2").Value = "N"
$ws.Range("C3").Value = "N"
$ws.Range("C4").Value = "N"

# --- New row 5: UnwatchArticleFromSearch ---
$ws.Range("A5").Value = "UnwatchArticleFromSearch"
$ws.Range("B5").Value = "To verify that user is able to unwatch a document from search results page"
$ws.Range("C5").Value = "Y"
$ws.Range("D5").Value = "SKIP"

# --- New row 6: UnwatchArticleFromRecordViewTest ---
$ws.Range("A6").Value = "UnwatchArticleFromRecordViewTest"
$ws.Range("B6").Value = "To verify that user is able to unwatch a document from document(Record View) page"
$ws.Range("C6").Value = "Y"
$ws.Range("D6").Value = "SKIP"

# --- New row 7: WatchlistArticleDocInfoTest ---
$ws.Range("A7").Value = "WatchlistArticleDocInfoTest"
$ws.Range("B7").Value = "To verify that the following fields are getting displayed for each document in watchlist page:`na)Times cited`nb)Comments`nc)Views"
$ws.Range("C7").Value = "Y"
$ws.Range("D7").Value = "PASS"

# Apply the same border/fill/wrap formatting used by the existing data rows to
# the freshly written cells (A/B/C/D for rows 5-7), then layer the distinct
# left/right-only border onto the Runmode cells (C5:C7) used to flag the new
# watchlist cases.
$ws.Range("A2:D4").Copy()
$ws.Range("A5:D7").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("C6").Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$ws.Range("C6").Borders.Item(10).LineStyle = 1  # xlEdgeRight

# Row 7 needs extra height to fit the three-line description
$ws.Rows.Item(7).RowHeight = 60

# Update the view: Excel originally had the window scrolled/selected on the
# old data; after adding rows the selection moves on
$ws.Range("B11").Select()

Write-Host "watchlist new test cases implementations applied"
